# The commit tidies up the "packet layout" worksheet (4th sheet):
#  - GPS lng/lat significance fields shrink from 8 bytes each to 4 bytes each
#  - Byte numbering (column A) is renumbered accordingly
#  - PD (photodiode) and EFM move up, and GPS gets a dedicated PD entry
#  - The last 8 rows (old bytes 27-34 / EFM, Cloud mean, Cloud st.dev, Rel hum)
#    are removed from the bottom, since the optional-section rows shift up
#  - The optional section (Cloud mean, Cloud st.dev, Rel hum) now starts
#    right after EFM, keeping the "OPTIONAL" marker in column D
#  - Final table covers rows 1-28 instead of rows 1-36

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Remove the now-unused trailing rows (old rows 29-36) first, shrinking the
# table from 36 data rows down to 28.
$ws.Rows("29:36").Delete() | Out-Null

# Rewrite the byte table from row 11 (byte 9) through row 28 (byte 30) to
# reflect the shrunk GPS fields and the shifted EFM / optional rows.
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "GPS lng"
$ws.Cells.Item(11,3).Value = 3

$ws.Cells.Item(12,1).Value = 14
$ws.Cells.Item(12,2).Value = "GPS lng"
$ws.Cells.Item(12,3).Value = 2

$ws.Cells.Item(13,1).Value = 15
$ws.Cells.Item(13,2).Value = "GPS lng"
$ws.Cells.Item(13,3).Value = 1

$ws.Cells.Item(14,1).Value = 16
$ws.Cells.Item(14,2).Value = "GPS lng"
$ws.Cells.Item(14,3).Value = 0

$ws.Cells.Item(15,1).Value = 17
$ws.Cells.Item(15,2).Value = "GPS lat"
$ws.Cells.Item(15,3).Value = 3

$ws.Cells.Item(16,1).Value = 18
$ws.Cells.Item(16,2).Value = "GPS lat"
$ws.Cells.Item(16,3).Value = 2

$ws.Cells.Item(17,1).Value = 19
$ws.Cells.Item(17,2).Value = "GPS lat"
$ws.Cells.Item(17,3).Value = 1

$ws.Cells.Item(18,1).Value = 20
$ws.Cells.Item(18,2).Value = "GPS lat"
$ws.Cells.Item(18,3).Value = 0

$ws.Cells.Item(19,1).Value = 21
$ws.Cells.Item(19,2).Value = "PD"
$ws.Cells.Item(19,3).Value = 1

$ws.Cells.Item(20,1).Value = 22
$ws.Cells.Item(20,2).Value = "PD"
$ws.Cells.Item(20,3).Value = 0

$ws.Cells.Item(21,1).Value = 23
$ws.Cells.Item(21,2).Value = "EFM"
$ws.Cells.Item(21,3).Value = 1
$ws.Cells.Item(21,4).Value = "OPTIONAL"

$ws.Cells.Item(22,1).Value = 24
$ws.Cells.Item(22,2).Value = "EFM"
$ws.Cells.Item(22,3).Value = 0
$ws.Cells.Item(22,4).Value = "OPTIONAL"

$ws.Cells.Item(23,1).Value = 25
$ws.Cells.Item(23,2).Value = "Cloud mean"
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = "OPTIONAL"

$ws.Cells.Item(24,1).Value = 26
$ws.Cells.Item(24,2).Value = "Cloud mean"
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = "OPTIONAL"

$ws.Cells.Item(25,1).Value = 27
$ws.Cells.Item(25,2).Value = "Cloud st.dev"
$ws.Cells.Item(25,3).Value = 1
$ws.Cells.Item(25,4).Value = "OPTIONAL"

$ws.Cells.Item(26,1).Value = 28
$ws.Cells.Item(26,2).Value = "Cloud st.dev"
$ws.Cells.Item(26,3).Value = 0
$ws.Cells.Item(26,4).Value = "OPTIONAL"

$ws.Cells.Item(27,1).Value = 29
$ws.Cells.Item(27,2).Value = "Rel hum"
$ws.Cells.Item(27,3).Value = 1
$ws.Cells.Item(27,4).Value = "OPTIONAL"

$ws.Cells.Item(28,1).Value = 30
$ws.Cells.Item(28,2).Value = "Rel hum"
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = "OPTIONAL"

# Restore the active-cell selection recorded for this sheet after the edit.
$ws.Range("F22").Select() | Out-Null
